$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.717406511306763
$ws.Range("B1").Value = 2.353511571884155
$ws.Range("C1").Value = 2.047903060913086
$ws.Range("D1").Value = 1.703803777694702
$ws.Range("E1").Value = 1.613725185394287
